$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColLetterToNum($letters) {
    $num = 0
    foreach ($ch in $letters.ToCharArray()) {
        $num = $num * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $num
}

function NumToColLetter($num) {
    $letters = ""
    while ($num -gt 0) {
        $rem = ($num - 1) % 26
        $letters = [char]([int][char]'A' + $rem) + $letters
        $num = [int](($num - $rem - 1) / 26)
    }
    return $letters
}

# The column being removed (culture_collection) and the last used column on the
# header/comment row before the edit.
$delCol = ColLetterToNum("V")
$lastCol = ColLetterToNum("CF")

# --- Step 1: re-home the cell comments that live on row 15 so they keep
# matching the field they describe once column V disappears. Comments are
# not moved automatically by a column delete, so shift them manually:
# collect the existing comment text for V..CF, wipe them, then re-create
# them one column to the left (skipping the one that belonged to V, which
# is the comment being dropped).
$colText = @{}
for ($col = $delCol; $col -le $lastCol; $col++) {
    $addr = "$(NumToColLetter($col))15"
    $rng = $ws.Range($addr)
    if ($rng.Comment -ne $null) {
        $colText[$col] = $rng.Comment.Text()
    }
}

for ($col = $delCol; $col -le $lastCol; $col++) {
    $addr = "$(NumToColLetter($col))15"
    $rng = $ws.Range($addr)
    if ($rng.Comment -ne $null) {
        $rng.Comment.Delete()
    }
}

for ($col = ($delCol + 1); $col -le $lastCol; $col++) {
    if ($colText.ContainsKey($col)) {
        $newAddr = "$(NumToColLetter($col - 1))15"
        $ws.Range($newAddr).AddComment($colText[$col])
    }
}

# --- Step 2: remove the culture_collection column itself (header + any
# data under it), shifting everything to its right one column to the left.
$ws.Range("V:V").Delete()
